$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.706.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.545.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.22%  "
$ws.Range("E13").Value = "  -1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.936.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.545.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.834"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.724.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0953"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("E24").Value = "  -1.55%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0793"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").Value = "  -2.69%  "
$ws.Range("E35").Value = "  -4.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.02%  "
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("E42").Value = "  +5.28%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  -1.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.991.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.790.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.17%  "
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.50%  "
